# Apply updated loading_percent values for rows 2-25, columns B-O (14 kV busbar columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each inner array holds the new B..O values (in column order) for one row.
$newValues = @(
    @(15.51149117908187, 11.31816739524666, 9.973725753208287, 0, 30.12093933366523, 29.46207227337908, 14.58530314740084, 22.90330466232728, 10.62414317297589, 0, 11.75369386101375, 15.81922631211973, 0, 22.27222467292891),
    @(14.92338489953312, 11.1195265447466, 9.985187642864499, 0, 30.2490390335542, 29.60978525125268, 14.64447923465788, 23.04189895437295, 10.64916100325633, 0, 11.75182347392452, 15.67643477404337, 0, 22.37732768848678),
    @(14.55018300177139, 10.99563146837196, 9.993416748024108, 0, 30.33506072136976, 29.7107304879151, 14.68329193825204, 23.13169962287188, 10.66534036450876, 0, 11.75186771082815, 15.58917950923295, 0, 22.44696958523119),
    @(14.39525485543039, 10.94470750471847, 9.997070284671601, 0, 30.37196401799001, 29.75443080534808, 14.69973187842016, 23.1694785936845, 10.67213994727439, 0, 11.75218688854104, 15.55375513026389, 0, 22.47663184946208),
    @(14.36936353037539, 10.93622677396244, 9.997695091725294, 0, 30.3782033252387, 29.76184170607703, 14.70249938029249, 23.17582335593648, 10.6732814939366, 0, 11.7522581115442, 15.54788175841359, 0, 22.48163465866283),
    @(14.54810482923856, 10.99494638792035, 9.993464805110928, 0, 30.33555093188021, 29.71130948009614, 14.6835111283781, 23.13220432498719, 10.66543122972369, 0, 11.75187079442604, 15.58870118976916, 0, 22.44736443007231),
    @(15.3113357642575, 11.25010542835324, 9.977430838889136, 0, 30.16357653086141, 29.5108689323072, 14.60519282977899, 22.9501169088775, 10.63259983617021, 0, 11.7528021807997, 15.76992127397374, 0, 22.30740285409144),
    @(16.70442027646896, 11.73308818728836, 9.955418966716804, 0, 29.88497071742011, 29.1996843089001, 14.47126412796695, 22.6302788967511, 10.57468351002037, 0, 11.7640354469137, 16.12740943806096, 0, 22.07355373833452),
    @(17.65603805414819, 12.07468889172002, 9.944964051428197, 0, 29.71625742466092, 29.0217108889122, 14.38483200964384, 22.41787651599918, 10.53603630527692, 0, 11.77794496588465, 16.38972189098439, 0, 21.92661896987399),
    @(18.07187358306349, 12.22669618548312, 9.941441766005859, 0, 29.647365659313, 28.95189426355816, 14.34810693896811, 22.32612877630951, 10.51929434377696, 0, 11.78548288029545, 16.5086383434459, 0, 21.86520051788125),
    @(18.22678529807224, 12.28372972706678, 9.940284667368077, 0, 29.62241176940973, 28.92707058838075, 14.33457289273097, 22.29208572444155, 10.51307461978842, 0, 11.78850965335472, 16.55358122384735, 0, 21.8427249188574),
    @(18.19353746273897, 12.27147065893798, 9.940526021294231, 0, 29.62773554002184, 28.93234482963922, 14.33747110453746, 22.29938639851818, 10.51440881585322, 0, 11.78785014759651, 16.54390635254451, 0, 21.8475306124044),
    @(18.0846700314759, 12.23139918428473, 9.941343033291092, 0, 29.64528993899845, 28.94981958493383, 14.34698601035103, 22.32331401513142, 10.51878023951483, 0, 11.7857284511867, 16.5123377543081, 0, 21.86333574573254),
    @(18.01764988672964, 12.20678426829772, 9.941866468939235, 0, 29.65619030633367, 28.96073395608203, 14.35286273381032, 22.33806147417646, 10.5214734867609, 0, 11.78445124054161, 16.49298873259165, 0, 21.87311878857148),
    @(17.62850979813664, 12.06468310765908, 9.945219008275989, 0, 29.72091811806928, 29.02649884616486, 14.38728426541563, 22.42397044188238, 10.53714726600096, 0, 11.77747655613562, 16.38193965617891, 0, 21.93074214042347),
    @(17.3853375252234, 11.97661079500429, 9.947591186051881, 0, 29.76264182192724, 29.06970640613309, 14.40906498552548, 22.47792057210122, 10.54697709057359, 0, 11.77350662383528, 16.31368858703945, 0, 21.96748286389423),
    @(17.24387065438736, 11.92563728678734, 9.949071766379767, 0, 29.78737946896513, 29.09560654108995, 14.42183677394699, 22.50941020277232, 10.55270993325608, 0, 11.77133726110182, 16.27439543250092, 0, 21.98912544077533),
    @(17.19570087247255, 11.90832543847108, 9.94959303751215, 0, 29.79588205444853, 29.10455555709101, 14.42620300509818, 22.52015091266237, 10.55466455776034, 0, 11.77062238946905, 16.26108599393671, 0, 21.99654081340723),
    @(17.41139007759656, 11.98601931930802, 9.947326645772387, 0, 29.75812371974068, 29.06499830618112, 14.40672112773184, 22.47212999880664, 10.54592251791502, 0, 11.77391743972512, 16.32095806572435, 0, 21.96351892104989),
    @(18.1167171583991, 12.2431838015657, 9.941098266800614, 0, 29.64010297978207, 28.94464292097758, 14.34418113008588, 22.31626691042366, 10.51749299132795, 0, 11.78634698205807, 16.52161284066838, 0, 21.85867214593064),
    @(18.56275186390265, 12.40815883860311, 9.938057337769036, 0, 29.56958081496439, 28.87539854237577, 14.30548143520642, 22.21847997471069, 10.49961240602781, 0, 11.7954740121019, 16.65222503768569, 0, 21.79470931480719),
    @(18.32609189438492, 12.32040486037802, 9.939586362467525, 0, 29.60661356500694, 28.91149036070665, 14.3259372748158, 22.27029789559286, 10.50909175765992, 0, 11.79051150514623, 16.5825726112094, 0, 21.82842938362489),
    @(17.39961689336808, 11.98176678000518, 9.947445880622267, 0, 29.76016401744542, 29.06712353870992, 14.40778000784451, 22.47474644358362, 10.54639903616456, 0, 11.77373135773882, 16.31767170258288, 0, 21.96530940050605),
    @(16.33962072398278, 11.60457105198225, 9.960367294137939, 0, 29.95403922317672, 29.27502438415561, 14.50539308338026, 22.71282959623653, 10.58966312356299, 0, 11.75999724982064, 16.03064061762145, 0, 22.13245625636194)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowVals = $newValues[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, 2 + $j).Value2 = $rowVals[$j]
    }
}
